$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

$box = $s.Shapes.AddTextbox(1, 470.79, 359.14, 205.16, 26.1)
$box.Name = "Shape 89"

$tf = $box.TextFrame
$tf.MarginLeft = 7.198818897637795
$tf.MarginRight = 7.198818897637795
$tf.MarginTop = 7.198818897637795
$tf.MarginBottom = 7.198818897637795
$tf.VerticalAnchor = 1
$tf.HorizontalAnchor = 0
$tf.AutoSize = 0

$tb = $tf.TextRange
$tb.Text = "(http://littlegreenriver.com/weblog/wp-content/uploads/mtv-diagram-730x1024.png)"
$tb.Font.Size = 8
$tb.Font.Color.RGB = 13421772

$box.Fill.Visible = $false
$box.Line.Visible = $false
